# Apply stock symbol list changes:
#  - remove ABNB (row 2)
#  - remove ARKX (originally row 12)
#  - rename BRK.B -> BRK-B
# Net effect: symbol list shrinks from 89 to 87 rows, and the BRK.B entry
# becomes BRK-B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename BRK.B to BRK-B first, while row numbers are still at their
# original positions (row 18).
$ws.Range("A18").Value = "BRK-B"

# Delete ARKX (row 12) before ABNB (row 2) so row numbers above it are
# unaffected by the later deletion of row 2.
$ws.Rows.Item(12).Delete()   # ARKX

# Delete ABNB (row 2)
$ws.Rows.Item(2).Delete()    # ABNB
